$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 252 (pushes existing rows 252:322 down to 253:323,
# carrying the date-formatted style on column D along with it).
$ws.Rows.Item(252).Insert()

# Populate the freshly inserted row 252 with the new price-observation record.
$ws.Range("A252").Value = 5
$ws.Range("B252").Value = "Macroferia Regional de Talca"
$ws.Range("C252").Value = "Maule"
$ws.Range("D252").Value = 44798
$ws.Range("E252").Value = 7
$ws.Range("F252").Value = "Fruta"
$ws.Range("G252").Value = 100102
$ws.Range("H252").Value = "Cítricos"
$ws.Range("I252").Value = 100102004
$ws.Range("J252").Value = "Mandarina"
$ws.Range("K252").Value = "Murcott"
$ws.Range("L252").Value = "Primera"
$ws.Range("M252").Value = 280
$ws.Range("N252").Value = 7000
$ws.Range("O252").Value = 7000
$ws.Range("P252").Value = 7000
$ws.Range("Q252").Value = "$/bandeja 10 kilos"
$ws.Range("R252").Value = "Provincia de Limarí"
$ws.Range("S252").Value = 700
$ws.Range("T252").Value = 10
